$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "F1 OSR 2011" (race results) - add race 13 ("13-ITA") results
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("F1 OSR 2011")

# Header / count row for the new race column N
$ws1.Range("N3").Value2 = 7
$ws1.Range("N4").Value2 = "13-ITA"

# Finishing positions for the pilots that scored points in race 13
$ws1.Range("N6").Value2  = 1   # Evgeny Egorenko
$ws1.Range("N8").Value2  = 2   # Ivan Egorov
$ws1.Range("N10").Value2 = 3   # Igor Peshkov
$ws1.Range("N13").Value2 = 4   # Maksim Prokoshun
$ws1.Range("N12").Value2 = 5   # Andrey Stasiukevich
$ws1.Range("N17").Value2 = 6   # Alexandr Zakirov
$ws1.Range("N30").Value2 = 7   # Ilya Ivashchenko

$excel.CalculateFull()

# ---------------------------------------------------------------------------
# Sheet "Лист2" (sorted rating table) - refresh with the new standings
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ratingRows = @(
    @(1,  0.9029999999999999,  "Evgeny Egorenko"),
    @(2,  0.42400000000000004, "Alexandr Zakirov"),
    @(3,  0.388,               "Igor Peshkov"),
    @(4,  0.35600000000000004, "Maksim Prokoshun"),
    @(5,  0.356,               "Andrey Vinokurov"),
    @(6,  0.323,               "Axsan Kalimulin"),
    @(7,  0.282,               "Dmitry Ivanov"),
    @(8,  0.25999999999999995, "Ivan Egorov"),
    @(9,  0.187,               "Andrey Stasiukevich"),
    @(10, 0.14200000000000002, "Petr Myakushin"),
    @(11, 0.11599999999999999, "Ilya Alexandrov"),
    @(12, 0.082,               "Alexey Makeev"),
    @(13, 0.08099999999999999, "Evgeny Peshkov"),
    @(14, 0.08,                "Robert Mardanov"),
    @(15, 0.07,                "Sergey Mazurin"),
    @(16, 0.063,               "Sergey Lozgachev"),
    @(17, 0.062,               "Nikita Kashin"),
    @(18, 0.05399999999999999, "Vadim Vrenere"),
    @(19, 0.05,                "Andrey Korneev"),
    @(20, 0.046,               "Artem Gusev"),
    @(21, 0.039,               "Yury Sbitnev"),
    @(22, 0.036,               "Roman Eazotov"),
    @(23, 0.03,                "Nikolay Kondratev"),
    @(24, 0.023,               "Ilya Ivashchenko"),
    @(25, 0.018,               "Sergey Bondarchuk"),
    @(26, 0.01,                "Sergey Protiv")
)

$r = 2
foreach ($row in $ratingRows) {
    $ws2.Cells.Item($r, 1).Value2 = $row[0]
    $ws2.Cells.Item($r, 2).Value2 = $row[1]
    $ws2.Cells.Item($r, 3).Value2 = $row[2]
    $r++
}

# Position column (A2:A27) is now shown in bold
$ws2.Range("A2:A27").Font.Bold = $true

# Update the saved selections to match the latest editing session
$ws2.Activate()
$ws2.Range("F12").Select()

$ws1.Activate()
$ws1.Range("V5:W30").Select()
$ws2.Activate()
